# Updates cryptos list prices/volumes (GitHub Actions data refresh).
# Numeric-looking Price values are written with a leading apostrophe so
# Excel stores them as literal text (matching the sheet's existing
# inline-string cells) instead of auto-converting to numbers; Style is
# reset to Normal afterwards so no stray number-format style lingers on
# the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '76.972.56'
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').Value = '2.959.18'
$ws.Range('E3').Value = '  +3.00%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''199.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.22%  '
$ws.Range('D6').Value = '''596.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('E8').Value = '  -0.38%  '
$ws.Range('E9').Value = '  +1.78%  '
$ws.Range('D10').Value = '2.961.89'
$ws.Range('E10').Value = '  +3.18%  '
$ws.Range('E11').Value = '  +12.87%  '
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('D13').Value = '3.505.67'
$ws.Range('E13').Value = '  +3.14%  '
$ws.Range('D14').Value = '''4.90'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').Value = '76.761.80'
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('E16').Value = '  +3.22%  '
$ws.Range('E17').Value = '  -0.29%  '
$ws.Range('D18').Value = '2.965.78'
$ws.Range('E18').Value = '  +3.03%  '
$ws.Range('D19').Value = '''13.49'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.73%  '
$ws.Range('D20').Value = '''8.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.81%  '
$ws.Range('D21').Value = '''374.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.96%  '
$ws.Range('E22').Value = '  +5.10%  '
$ws.Range('E23').Value = '  -3.55%  '
$ws.Range('D24').Value = '''72.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.78%  '
$ws.Range('D25').Value = '3.116.26'
$ws.Range('E25').Value = '  +2.91%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = '''4.28'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.42%  '
$ws.Range('D28').Value = '''9.69'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('E29').Value = '  +2.85%  '
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('D31').Value = '''8.40'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +8.80%  '
$ws.Range('E32').Value = '  -1.74%  '
$ws.Range('D33').Value = '''498.89'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.08%  '
$ws.Range('E34').Value = '  +1.40%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = '''166.28'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.55%  '
$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').Value = '''0.400'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +15.91%  '
$ws.Range('B38').Value = 'Cronos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D38').Value = '''0.113'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +21.96%  '
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('E40').Value = '  +1.46%  '
$ws.Range('D41').Value = '''0.111'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.85%  '
$ws.Range('D43').Value = '''180.13'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.42%  '
$ws.Range('E45').Value = '  -1.48%  '
$ws.Range('D46').Value = '''40.15'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.19%  '
$ws.Range('E47').Value = '  -3.26%  '
$ws.Range('D48').Value = '''0.590'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.84%  '
$ws.Range('E49').Value = '  +4.14%  '
$ws.Range('D50').Value = '''2.30'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.54%  '
$ws.Range('E51').Value = '  +4.99%  '
